$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15
$ws.Range("A15").Value = 98151268
$ws.Range("B15").Value = 73631
$ws.Range("D15").Value = "LC"
$ws.Range("E15").Value = 6426
$ws.Range("F15").Value = "Kattfotslav"
$ws.Range("G15").Value = "Felipes leucopellaeus"
$ws.Range("H15").Value = "(Ach.) Frisch & G.Thor"
$ws.Range("P15").Value = "Bråtan, Ång"
$ws.Range("Q15").Value = 664289.4311076899
$ws.Range("R15").Value = 6968026.761211542

# Row 16
$ws.Range("A16").Value = 98150989
$ws.Range("B16").Value = 77541
$ws.Range("D16").Value = "NT"
$ws.Range("E16").Value = 185
$ws.Range("F16").Value = "Violettgrå tagellav"
$ws.Range("G16").Value = "Bryoria nadvornikiana"
$ws.Range("H16").Value = "(Gyeln.) Brodo & D.Hawksw."
$ws.Range("P16").Value = "Bråtan, Ång"
$ws.Range("Q16").Value = 664384.6636208369
$ws.Range("R16").Value = 6968003.128549194

# Row 17
$ws.Range("A17").Value = 98151014
$ws.Range("B17").Value = 89732
$ws.Range("D17").Value = "VU"
$ws.Range("E17").Value = 2062
$ws.Range("F17").Value = "Ulltickeporing"
$ws.Range("G17").Value = "Skeletocutis brevispora"
$ws.Range("H17").Value = "Niemelä"
$ws.Range("P17").Value = "Bråtan, Ång"
$ws.Range("Q17").Value = 664290.9406681373
$ws.Range("R17").Value = 6968060.737959783

# Row 18
$ws.Range("A18").Value = 98151080
$ws.Range("B18").Value = 89392
$ws.Range("D18").Value = "NT"
$ws.Range("E18").Value = 1202
$ws.Range("F18").Value = "Ullticka"
$ws.Range("G18").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H18").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("P18").Value = "Bråtan, Ång"
$ws.Range("Q18").Value = 664386.6082164289
$ws.Range("R18").Value = 6968010.097956684

# Row 19
$ws.Range("A19").Value = 98151357
$ws.Range("B19").Value = 78479
$ws.Range("D19").Value = "VU"
$ws.Range("E19").Value = 392
$ws.Range("F19").Value = "Aspgelélav"
$ws.Range("G19").Value = "Collema subnigrescens"
$ws.Range("H19").Value = "Degel."
$ws.Range("P19").Value = "Bråtan, Ång"
$ws.Range("Q19").Value = 664139.5142268437
$ws.Range("R19").Value = 6968033.915686903

# Row 20
$ws.Range("A20").Value = 98151084
$ws.Range("B20").Value = 89392
$ws.Range("D20").Value = "NT"
$ws.Range("E20").Value = 1202
$ws.Range("F20").Value = "Ullticka"
$ws.Range("G20").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H20").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("P20").Value = "Bråtan, Ång"
$ws.Range("Q20").Value = 664311.4867898196
$ws.Range("R20").Value = 6968099.332879714

# Row 21
$ws.Range("A21").Value = 98151131
$ws.Range("B21").Value = 78503
$ws.Range("D21").Value = "LC"
$ws.Range("E21").Value = 6456
$ws.Range("F21").Value = "Skinnlav"
$ws.Range("G21").Value = "Leptogium saturninum"
$ws.Range("H21").Value = "(Dicks.) Nyl."
$ws.Range("P21").Value = "Bråtan, Ång"
$ws.Range("Q21").Value = 664135.3270324104
$ws.Range("R21").Value = 6968035.080522615

# Row 22
$ws.Range("A22").Value = 98151083
$ws.Range("B22").Value = 89392
$ws.Range("D22").Value = "NT"
$ws.Range("E22").Value = 1202
$ws.Range("F22").Value = "Ullticka"
$ws.Range("G22").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H22").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("P22").Value = "Bråtan, Ång"
$ws.Range("Q22").Value = 664289.8194318935
$ws.Range("R22").Value = 6968064.804886819

# Row 23
$ws.Range("A23").Value = 98151081
$ws.Range("B23").Value = 89392
$ws.Range("D23").Value = "NT"
$ws.Range("E23").Value = 1202
$ws.Range("F23").Value = "Ullticka"
$ws.Range("G23").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H23").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("P23").Value = "Bråtan, Ång"
$ws.Range("Q23").Value = 664390.9109275232
$ws.Range("R23").Value = 6967988.323297421

# Row 24
$ws.Range("A24").Value = 98151176
$ws.Range("B24").Value = 89673
$ws.Range("D24").Value = "NT"
$ws.Range("E24").Value = 658
$ws.Range("F24").Value = "Rosenticka"
$ws.Range("G24").Value = "Rhodofomes roseus"
$ws.Range("H24").Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Range("P24").Value = "Bråtan, Ång"
$ws.Range("Q24").Value = 664312.310446391
$ws.Range("R24").Value = 6968101.206605153

# Row 25
$ws.Range("A25").Value = 98151145
$ws.Range("B25").Value = 89403
$ws.Range("D25").Value = "LC"
$ws.Range("E25").Value = 1205
$ws.Range("F25").Value = "Stor aspticka"
$ws.Range("G25").Value = "Phellinus populicola"
$ws.Range("H25").Value = "Niemelä"
$ws.Range("P25").Value = "Bråtan, Ång"
$ws.Range("Q25").Value = 664122.421552405
$ws.Range("R25").Value = 6968017.942315456

# Row 26
$ws.Range("A26").Value = 98151307
$ws.Range("B26").Value = 81236
$ws.Range("D26").Value = "NT"
$ws.Range("E26").Value = 1312
$ws.Range("F26").Value = "Gammelgransskål"
$ws.Range("G26").Value = "Pseudographis pinicola"
$ws.Range("H26").Value = "(Nyl.) Rehm"
$ws.Range("P26").Value = "Kråkan, Ång"
$ws.Range("Q26").Value = 664448.9567745975
$ws.Range("R26").Value = 6968056.745323926

# Row 27
$ws.Range("A27").Value = 98151269
$ws.Range("B27").Value = 73631
$ws.Range("D27").Value = "LC"
$ws.Range("E27").Value = 6426
$ws.Range("F27").Value = "Kattfotslav"
$ws.Range("G27").Value = "Felipes leucopellaeus"
$ws.Range("H27").Value = "(Ach.) Frisch & G.Thor"
$ws.Range("P27").Value = "Kråkan, Ång"
$ws.Range("Q27").Value = 664416.9957306259
$ws.Range("R27").Value = 6967989.630702551

# Row 28
$ws.Range("A28").Value = 98151342
$ws.Range("B28").Value = 77506
$ws.Range("D28").Value = "NT"
$ws.Range("E28").Value = 6425
$ws.Range("F28").Value = "Garnlav"
$ws.Range("G28").Value = "Alectoria sarmentosa"
$ws.Range("H28").Value = "(Ach.) Ach."
$ws.Range("P28").Value = "Kråkan, Ång"
$ws.Range("Q28").Value = 664414.5233459049
$ws.Range("R28").Value = 6968011.497083161

# Row 29
$ws.Range("A29").Value = 98151305
$ws.Range("B29").Value = 81236
$ws.Range("D29").Value = "NT"
$ws.Range("E29").Value = 1312
$ws.Range("F29").Value = "Gammelgransskål"
$ws.Range("G29").Value = "Pseudographis pinicola"
$ws.Range("H29").Value = "(Nyl.) Rehm"
$ws.Range("P29").Value = "Kråkan, Ång"
$ws.Range("Q29").Value = 664420.7254760786
$ws.Range("R29").Value = 6967988.443260053

# Row 30
$ws.Range("A30").Value = 98151275
$ws.Range("B30").Value = 89406
$ws.Range("D30").Value = "NT"
$ws.Range("E30").Value = 1204
$ws.Range("F30").Value = "Gränsticka"
$ws.Range("G30").Value = "Phellopilus nigrolimitatus"
$ws.Range("H30").Value = "(Romell) Niemelä, T.Wagner & M.Fisch."
$ws.Range("P30").Value = "Kråkan, Ång"
$ws.Range("Q30").Value = 664428.9819053004
$ws.Range("R30").Value = 6968052.537048084

# Row 31
$ws.Range("A31").Value = 98151283
$ws.Range("B31").Value = 89410
$ws.Range("D31").Value = "NT"
$ws.Range("E31").Value = 5432
$ws.Range("F31").Value = "Granticka"
$ws.Range("G31").Value = "Porodaedalea chrysoloma"
$ws.Range("H31").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("P31").Value = "Kråkan, Ång"
$ws.Range("Q31").Value = 664467.6032586561
$ws.Range("R31").Value = 6968078.295775652

# Row 32
$ws.Range("A32").Value = 98151358
$ws.Range("B32").Value = 93148
$ws.Range("D32").Value = "VU"
$ws.Range("E32").Value = 1079
$ws.Range("F32").Value = "Aspfjädermossa"
$ws.Range("G32").Value = "Neckera pennata"
$ws.Range("H32").Value = "Hedw."
$ws.Range("P32").Value = "Stormyran, Ång"
$ws.Range("Q32").Value = 663923.3061348161
$ws.Range("R32").Value = 6968220.089407885

# Row 33
$ws.Range("A33").Value = 98151130
$ws.Range("B33").Value = 78503
$ws.Range("D33").Value = "LC"
$ws.Range("E33").Value = 6456
$ws.Range("F33").Value = "Skinnlav"
$ws.Range("G33").Value = "Leptogium saturninum"
$ws.Range("H33").Value = "(Dicks.) Nyl."
$ws.Range("P33").Value = "Bråtan, Ång"
$ws.Range("Q33").Value = 663974.0132528287
$ws.Range("R33").Value = 6968297.750628498

# Row 34
$ws.Range("A34").Value = 98151356
$ws.Range("B34").Value = 78479
$ws.Range("D34").Value = "VU"
$ws.Range("E34").Value = 392
$ws.Range("F34").Value = "Aspgelélav"
$ws.Range("G34").Value = "Collema subnigrescens"
$ws.Range("H34").Value = "Degel."
$ws.Range("P34").Value = "Bråtan, Ång"
$ws.Range("Q34").Value = 663975.9349722316
$ws.Range("R34").Value = 6968296.014373922

# Row 35
$ws.Range("A35").Value = 98151144
$ws.Range("B35").Value = 89403
$ws.Range("D35").Value = "LC"
$ws.Range("E35").Value = 1205
$ws.Range("F35").Value = "Stor aspticka"
$ws.Range("G35").Value = "Phellinus populicola"
$ws.Range("H35").Value = "Niemelä"
$ws.Range("P35").Value = "Bråtan, Ång"
$ws.Range("Q35").Value = 663973.143774787
$ws.Range("R35").Value = 6968296.79099117
